$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header area -----------------------------------------------------
# Date (G4) -- new timestamp (Excel serial date/time, style s="59" unchanged)
$ws.Range("G4").Value = 44084.83020343081

# Customer name (G7)
$ws.Range("G7").Value = "Nueva propuesta"

# Address comments box (F10)
$ws.Range("F10").Value = "vamos"

# --- Offline storage section ------------------------------------------
# Digital (GB): Qty 450 -> 240, total 6750 -> 3600
$ws.Range("F19").Value = 240
$ws.Range("H19").Value = 3600

# Visual (pages) row cleared entirely (Pages/frame qty + price cells)
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = ""

# --- Registration fee (row 24) -----------------------------------------
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 200
$ws.Range("H24").Value = 200

# --- AWA contribution (row 25) ------------------------------------------
$ws.Range("E25").Value = "public"
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 500
$ws.Range("H25").Value = 500

# --- Management fee (row 26) --------------------------------------------
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 60
$ws.Range("H26").Value = 300

# --- Storage per reel / per period (row 27) -----------------------------
# E27 is numeric-looking text ("5"): enter it quote-prefixed so it is
# stored as text, then re-apply the original number format (copied from
# the sibling E22 cell which shares the same style index) so the cell
# keeps its original style id instead of gaining a new quote-prefixed one.
$ws.Range("E27").Value = "'5"
$ws.Range("E22").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 60
$ws.Range("H27").Value = 600

# --- Professional services (row 28) -------------------------------------
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 900
$ws.Range("H28").Value = 900

# --- piqlReader (row 29) --------------------------------------------------
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 79900
$ws.Range("H29").Value = 79900

# --- Installation and training (row 30) -----------------------------------
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 3000
$ws.Range("H30").Value = 3000

# --- Service agreement (row 31) -------------------------------------------
$ws.Range("E31").Value = "gold"
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 2500
$ws.Range("H31").Value = 2500

# --- Shipment cost (row 32) -----------------------------------------------
$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 20
$ws.Range("H32").Value = 40

# --- Totals ----------------------------------------------------------------
$ws.Range("H33").Value = 93000
$ws.Range("H34").Value = 2500
